$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 21924816
$ws.Range("B3").Value = "Bendel Dan"
$ws.Range("C3").Value = "Testing independence using random projections"
$ws.Range("E3").Value = "matlab"

$ws.Range("A4").Value = 38000014
$ws.Range("B4").Value = "Fainblat Ido"
$ws.Range("C4").Value = "Testing independence using random projections"
$ws.Range("E4").Value = "matlab"
